$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 137, shifting existing rows 137-200 down to 138-201.
$ws.Rows(137).Insert()

# Populate the newly inserted row 137 with the new weekly entry.
$ws.Range("A137").Value = 3
$ws.Range("B137").Value = "Femacal de La Calera"
$ws.Range("C137").Value = "Coquimbo"
$ws.Range("D137").Value = 44466
$ws.Range("E137").Value = 5
$ws.Range("F137").Value = 100112043
$ws.Range("G137").Value = "Pepino ensalada"
$ws.Range("H137").Value = "Sin especificar"
$ws.Range("I137").Value = "Primera"
$ws.Range("J137").Value = 130
$ws.Range("K137").Value = 11500
$ws.Range("L137").Value = 12000
$ws.Range("M137").Value = 11731
$ws.Range("N137").Value = "$/caja 70 unidades"
$ws.Range("O137").Value = "Región de Arica y Parinacota"
$ws.Range("P137").Value = 168
$ws.Range("Q137").Value = 70
$ws.Range("R137").Value = "Hortaliza"
